$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add 6 new "France" document-type rows (272-277) below the existing data,
# replacing the previously-blank trailing row 272.
# ---------------------------------------------------------------------------

$cols = @("A","B","C","D","E","H","I","J","L","M","N")

# 1) Copy the cell formatting (styles) from the template row (271) onto each
#    of the 6 new rows, cell by cell, so that only the same columns that are
#    populated in row 271 receive an explicit style (this mirrors the way
#    the sheet lays out its data and avoids stray formatting on columns
#    F/G/K, which stay empty).
for ($r = 272; $r -le 277; $r++) {
    foreach ($col in $cols) {
        $ws.Range($col + "271").Copy() | Out-Null
        $ws.Range($col + $r).PasteSpecial(-4122, -4142, $false, $false) | Out-Null
    }
}
$excel.CutCopyMode = $false

# 2) Give the new rows the same (taller, wrapped-text) row height used by
#    their neighbours.
$ws.Range("A272:N277").RowHeight = 30

# 3) Populate the common columns shared by every new row.
for ($r = 272; $r -le 277; $r++) {
    $ws.Range("B" + $r).Value = "busdox-docid-qns"
    $ws.Range("D" + $r).Value = "9.0"
    $ws.Range("E" + $r).Value = "active"
    $ws.Range("H" + $r).Value = "TICC-363"
    $ws.Range("I" + $r).Value = $false
    $ws.Range("J" + $r).Formula = "=TRUE"
    $ws.Range("L" + $r).Value = "POAC-France"
    $ws.Range("N" + $r).Value = "cenbii-procid-ubl::urn:peppol:france:billing:regulated`r`ncenbii-procid-ubl::urn:peppol:france:billing:non-regulated"
}

# 4) Populate the row-specific columns (document type identifier in column C
#    and category/name in columns A and M).
$ws.Range("C272").Value = "urn:peppol:doctype:pdf+xml##urn:cen.eu:en16931:2017#conformant#urn:peppol:france:billing:Factur-X:1.0::D22B"
$ws.Range("C273").Value = "urn:oasis:names:specification:ubl:schema:xsd:Invoice-2::Invoice##urn:cen.eu:en16931:2017#compliant#urn:peppol:france:billing:cius:1.0::2.1"
$ws.Range("C274").Value = "urn:oasis:names:specification:ubl:schema:xsd:Invoice-2::Invoice##urn:cen.eu:en16931:2017#conformant#urn:peppol:france:billing:extended:1.0::2.1"

$ws.Range("A273").Value = "France UBL Invoice CIUS"
$ws.Range("A274").Value = "France UBL Invoice Extension"
$ws.Range("A275").Value = "France CII Invoice CIUS"
$ws.Range("A276").Value = "France CII Invoice Extension"

$ws.Range("C275").Value = "urn:un:unece:uncefact:data:standard:CrossIndustryInvoice:100::CrossIndustryInvoice##urn:cen.eu:en16931:2017#compliant#urn:peppol:france:billing:cius:1.0::D16B"
$ws.Range("C276").Value = "urn:un:unece:uncefact:data:standard:CrossIndustryInvoice:100::CrossIndustryInvoice##urn:cen.eu:en16931:2017#conformant#urn:peppol:france:billing:extended:1.0::D22B"
$ws.Range("C277").Value = "urn:un:unece:uncefact:data:standard:CrossDomainAcknowledgementAndResponse:100::CrossDomainAcknowledgementAndResponse##urn:peppol:france:billing:cdv:1.0::D22B"

$ws.Range("A272").Value = "France Factur-X"
$ws.Range("A277").Value = "France CDAR"

for ($r = 272; $r -le 276; $r++) {
    $ws.Range("M" + $r).Value = "Invoice"
}
$ws.Range("M277").Value = "Invoice Response"

# 5) Match the workbook's on-screen selection to the newly entered data.
$ws.Range("A274").Select() | Out-Null
